# The deck shipped with its two theme parts "crossed": the theme that is
# actually wired to the slide master / presentation (the one driving every
# slide's look) held the "Integral" palette, while the unused theme that's
# only wired to the notes master held the plain "Office Theme" palette.
#
# The authored change swaps the two themes' contents so the palette that
# actually paints the deck becomes the stock "Office Theme" colors
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) instead of "Integral" - the
# font scheme and format (fill/line/effect) scheme are identical between
# the two themes, so the only observable difference is the 12 color-scheme
# entries. We rewrite them through the live presentation's theme color
# scheme, which is backed by the same theme part the slide master/
# presentation point at.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$colors = $slide.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink.
# Target values are the stock Office Theme RGBs (as 0xBBGGRR long values,
# the same packing PowerPoint's RGB property uses).
$officeThemeRGB = @(
    0,          # Dark1   000000
    16777215,   # Light1  FFFFFF
    6968388,    # Dark2   44546A
    15132391,   # Light2  E7E6E6
    13998939,   # Accent1 5B9BD5
    3243501,    # Accent2 ED7D31
    10855845,   # Accent3 A5A5A5
    49407,      # Accent4 FFC000
    12874308,   # Accent5 4472C4
    4697456,    # Accent6 70AD47
    12673797,   # Hyperlink       0563C1
    7491477     # FollowedHyperlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
